$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column F ("Is Significant"), which shifts it to H
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(6).Insert()

# New header cells
$ws.Range("F1").Value = "Observed"
$ws.Range("G1").Value = "Expected"

# Observed / Expected contingency-table strings per feature row
$observed = @(
    '[768 188] ; [ 76 142]',
    '[734 245] ; [107 118]',
    '[769 199] ; [133  92]',
    '[780 188] ; [141  81]',
    '[472 485] ; [ 65 156]',
    '[576 371] ; [ 93 129]',
    '[831 148] ; [159  66]',
    '[699 269] ; [123 100]',
    '[393 577] ; [ 51 171]',
    '[540 426] ; [ 85 139]',
    '[895  75] ; [182  41]',
    '[893  78] ; [182  41]',
    '[603 352] ; [103 119]',
    '[452 518] ; [ 70 154]',
    '[462 506] ; [ 72 150]',
    '[423 502] ; [ 65 146]',
    '[303 635] ; [ 41 176]',
    '[485 471] ; [ 81 141]',
    '[505 435] ; [ 89 130]',
    '[687 291] ; [131  94]',
    '[724 251] ; [141  82]'
)

$expected = @(
    '[687.27768313 268.72231687] ; [156.72231687  61.27768313]',
    '[683.83637874 295.16362126] ; [157.16362126  67.83637874]',
    '[731.88264878 236.11735122] ; [170.11735122  54.88264878]',
    '[749.18319328 218.81680672] ; [171.81680672  50.18319328]',
    '[436.25551783 520.74448217] ; [100.74448217 120.25551783]',
    '[541.95295124 405.04704876] ; [127.04704876  94.95295124]',
    '[804.99169435 174.00830565] ; [185.00830565  39.99169435]',
    '[668.0906801 299.9093199] ; [153.9093199  69.0906801]',
    '[361.30872483 608.69127517] ; [ 82.69127517 139.30872483]',
    '[507.35294118 458.64705882] ; [117.64705882 106.35294118]',
    '[875.68315172  94.31684828] ; [201.31684828  21.68315172]',
    '[874.22529313  96.77470687] ; [200.77470687  22.22529313]',
    '[572.83772302 382.16227698] ; [133.16227698  88.83772302]',
    '[424.07035176 545.92964824] ; [ 97.92964824 126.07035176]',
    '[434.37983193 533.62016807] ; [ 99.62016807 122.37983193]',
    '[397.35915493 527.64084507] ; [ 90.64084507 120.35915493]',
    '[279.36969697 658.63030303] ; [ 64.63030303 152.36969697]',
    '[459.3344652 496.6655348] ; [106.6655348 115.3344652]',
    '[481.76013805 458.23986195] ; [112.23986195 106.76013805]',
    '[665.0074813 312.9925187] ; [152.9925187  72.0074813]',
    '[703.98580968 271.01419032] ; [161.01419032  61.98580968]'
)

for ($i = 0; $i -lt 21; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $observed[$i]
    $ws.Cells.Item($row, 7).Value = $expected[$i]
}

Write-Output "Done"
